# Applies the "Trade #15 closed at 2026-02-17 20:52:33" update:
#  - Summary metrics refreshed
#  - Strategy Status for MarketMaking refreshed
#  - Trade #43 (row 44 in "All Trades", row 11 in "MarketMaking") closed early
#  - New trade #76 opened (new row 77 in "All Trades", new row 44 in "MarketMaking")

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell without Excel's automatic
# date/time re-interpretation, and without leaving a residual number-format
# style attached to the cell afterwards.
function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value2 = 1400.4
$summary.Range("B4").Value2 = 0.19
$summary.Range("B5").Value2 = 0.09
$summary.Range("B6").Value2 = 43
$summary.Range("B8").Value2 = 18
$summary.Range("B9").Value2 = 44.19

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value2 = 100.4
$status.Range("D5").Value2 = 10
$status.Range("E5").Value2 = 0.08
$status.Range("F5").Value2 = 0.4
$status.Range("G5").Value2 = 50

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #43 -> row 44 : close it out early
$allTrades.Cells.Item(44, 7).Value2 = 0.135764
$allTrades.Cells.Item(44, 8).Value2 = "CLOSED"
$allTrades.Cells.Item(44, 9).Value2 = -20.1387
$allTrades.Cells.Item(44, 10).Value2 = -0.03
$allTrades.Cells.Item(44, 11).Value2 = 100.4
Set-TextCell $allTrades.Cells.Item(44, 12) "early_exit"
$allTrades.Cells.Item(44, 13).Value2 = 0.13

# New trade #76 -> new row 77
$allTrades.Cells.Item(77, 1).Value2 = 76
Set-TextCell $allTrades.Cells.Item(77, 2) "2026-02-17"
Set-TextCell $allTrades.Cells.Item(77, 3) "20:52:26"
Set-TextCell $allTrades.Cells.Item(77, 4) "MarketMaking"
Set-TextCell $allTrades.Cells.Item(77, 5) "UP"
$allTrades.Cells.Item(77, 6).Value2 = 0.17
Set-TextCell $allTrades.Cells.Item(77, 8) "OPEN"
$allTrades.Cells.Item(77, 9).Value2 = 0
$allTrades.Cells.Item(77, 10).Value2 = 0
$allTrades.Cells.Item(77, 11).Value2 = 100.436797675607
$allTrades.Cells.Item(77, 13).Value2 = 0
$allTrades.Cells.Item(77, 14).Value2 = 0
$allTrades.Cells.Item(77, 15).Value2 = 0
$allTrades.Cells.Item(77, 16).Value2 = 0.6
Set-TextCell $allTrades.Cells.Item(77, 17) "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

# Trade #43 -> row 11 : close it out early
$marketMaking.Cells.Item(11, 7).Value2 = 0.135764
$marketMaking.Cells.Item(11, 8).Value2 = "CLOSED"
$marketMaking.Cells.Item(11, 9).Value2 = -20.1387
$marketMaking.Cells.Item(11, 10).Value2 = -0.03
$marketMaking.Cells.Item(11, 11).Value2 = 100.4
Set-TextCell $marketMaking.Cells.Item(11, 16) "early_exit"
$marketMaking.Cells.Item(11, 17).Value2 = 0.13

# New trade #76 -> new row 44
$marketMaking.Cells.Item(44, 1).Value2 = 76
Set-TextCell $marketMaking.Cells.Item(44, 2) "2026-02-17"
Set-TextCell $marketMaking.Cells.Item(44, 3) "20:52:26"
Set-TextCell $marketMaking.Cells.Item(44, 4) "MarketMaking"
Set-TextCell $marketMaking.Cells.Item(44, 5) "UP"
$marketMaking.Cells.Item(44, 6).Value2 = 0.17
Set-TextCell $marketMaking.Cells.Item(44, 8) "OPEN"
$marketMaking.Cells.Item(44, 9).Value2 = 0
$marketMaking.Cells.Item(44, 10).Value2 = 0
$marketMaking.Cells.Item(44, 11).Value2 = 100.436797675607
$marketMaking.Cells.Item(44, 12).Value2 = 0
$marketMaking.Cells.Item(44, 13).Value2 = 0
$marketMaking.Cells.Item(44, 14).Value2 = 0.6
Set-TextCell $marketMaking.Cells.Item(44, 15) "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item(44, 17).Value2 = 0
